$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$wsMeta.Range("B15").Value = "4.0.1"

# --- Elements sheet updates ---
$wsElem = $wb.Worksheets.Item("Elements")

# Extension row: constraint text simplified (no longer mentions Parameters resource)
$wsElem.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + "`n" + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Extension.id row: type changes from id to string
$wsElem.Range("K3").Value = "string`n"

# Extension.value[x] row: FHIR R4B -> R4 link
$wsElem.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
